# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the Seraph_Profits workbook sheets as described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 188.9  # was 306.07693
$ws.Range("I6").Value = 117.375  # was 97.7
$ws.Range("J6").Value = 475  # was 1000.6667
$ws.Range("K6").Value = 352.125  # was 293.1
$ws.Range("L6").Value = 1425  # was 3002.0001
$ws.Range("M6").Value = -240.125  # was -181.1
$ws.Range("N6").Value = -1649  # was -3226.0001

$ws.Range("H40").Value = 1989.9  # was 1988.8889
$ws.Range("I40").Value = 1989.9  # was 1988.8889
$ws.Range("K40").Value = 1989.9  # was 1988.8889
$ws.Range("M40").Value = -1814.9  # was -1813.8889

$ws.Range("H80").Value = 1248.1578  # was 4822.4736
$ws.Range("I80").Value = 830.6667  # was 570.6
$ws.Range("J80").Value = 1326.4375  # was 6341
$ws.Range("K80").Value = 2492.0001  # was 1711.8
$ws.Range("L80").Value = 3979.3125  # was 19023
$ws.Range("M80").Value = -1494.0001  # was -713.8000000000002
$ws.Range("N80").Value = -5975.3125  # was -21019

$ws.Range("H83").Value = 1248.1578  # was 4822.4736
$ws.Range("I83").Value = 830.6667  # was 570.6
$ws.Range("J83").Value = 1326.4375  # was 6341
$ws.Range("K83").Value = 7476.0003  # was 5135.400000000001
$ws.Range("L83").Value = 11937.9375  # was 57069
$ws.Range("M83").Value = -2484.0003  # was -143.4000000000005
$ws.Range("N83").Value = -21921.9375  # was -67053

$ws.Range("H107").Value = 55556960  # was 62501570
$ws.Range("I107").Value = 58825004  # was 62501570
$ws.Range("J107").Value = 200  # was 0
$ws.Range("K107").Value = 58825004  # was 62501570
$ws.Range("L107").Value = 200  # was 0
$ws.Range("M107").Value = -58823084  # was -62499650
$ws.Range("N107").Value = -4040  # new cell

$ws.Range("H120").Value = 0  # was 50000
$ws.Range("J120").Value = 0  # was 50000
$ws.Range("L120").Value = 0  # was 50000
$ws.Range("N120").ClearContents()  # was -59676

$ws.Range("H137").Value = 2271.2666  # was 2333.2559
$ws.Range("I137").Value = 2868.3333  # was 3071.4736
$ws.Range("K137").Value = 8604.999899999999  # was 9214.4208
$ws.Range("M137").Value = -6054.999899999999  # was -6664.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1489.3334  # was 2480.8
$ws.Range("I16").Value = 1550.625  # was 2851.25
$ws.Range("K16").Value = 1550.625  # was 2851.25
$ws.Range("M16").Value = -1263.625  # was -2564.25

$ws.Range("H32").Value = 390294.44  # was 422788.62
$ws.Range("I32").Value = 5092.4443  # was 5683.5
$ws.Range("K32").Value = 5092.4443  # was 5683.5
$ws.Range("M32").Value = -4805.4443  # was -5396.5

$ws.Range("H44").Value = 29883.4  # was 27059.428
$ws.Range("J44").Value = 37249.25  # was 31499.334
$ws.Range("L44").Value = 37249.25  # was 31499.334
$ws.Range("N44").Value = -38225.25  # was -32475.334

$ws.Range("H55").Value = 16948.715  # was 17626.666
$ws.Range("J55").Value = 25999.5  # was 23999.5
$ws.Range("L55").Value = 25999.5  # was 23999.5
$ws.Range("N55").Value = -26629.5  # was -24629.5

$ws.Range("H110").Value = 7937793  # was 10102511
$ws.Range("I110").Value = 13889513  # was 22222924
$ws.Range("K110").Value = 13889513  # was 22222924
$ws.Range("M110").Value = -13887468  # was -22220879

$ws.Range("H125").Value = 69998  # was 69996
$ws.Range("J125").Value = 69998  # was 69996
$ws.Range("L125").Value = 69998  # was 69996
$ws.Range("N125").Value = -79838  # was -79836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1839  # was 1748.7273
$ws.Range("J12").Value = 1697.75  # was 1527.4
$ws.Range("L12").Value = 1697.75  # was 1527.4
$ws.Range("N12").Value = -2033.75  # was -1863.4

$ws.Range("H22").Value = 200388.2  # was 433.66666
$ws.Range("I22").Value = 547  # was 501
$ws.Range("J22").Value = 500150  # was 400
$ws.Range("K22").Value = 547  # was 501
$ws.Range("L22").Value = 500150  # was 400
$ws.Range("M22").Value = -374  # was -328
$ws.Range("N22").Value = -500496  # was -746

$ws.Range("H94").Value = 1221.25  # was 1328.4
$ws.Range("I94").Value = 1205.75  # was 1326.4286
$ws.Range("J94").Value = 1252.25  # was 1333
$ws.Range("K94").Value = 1205.75  # was 1326.4286
$ws.Range("L94").Value = 1252.25  # was 1333
$ws.Range("M94").Value = -754.75  # was -875.4286
$ws.Range("N94").Value = -2154.25  # was -2235

$ws.Range("H105").Value = 10418117  # was 9260657
$ws.Range("I105").Value = 11906119  # was 10417977
$ws.Range("K105").Value = 11906119  # was 10417977
$ws.Range("M105").Value = -11904372  # was -10416230

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 755.1667  # was 672.0714
$ws.Range("J2").Value = 425  # was 341.16666
$ws.Range("L2").Value = 425  # was 341.16666
$ws.Range("N2").Value = -651  # was -567.16666

$ws.Range("H31").Value = 4219.6787  # was 4402
$ws.Range("I31").Value = 4207.6  # was 4415.357
$ws.Range("J31").Value = 4233.615  # was 4386.4165
$ws.Range("K31").Value = 4207.6  # was 4415.357
$ws.Range("L31").Value = 4233.615  # was 4386.4165
$ws.Range("M31").Value = -3912.6  # was -4120.357
$ws.Range("N31").Value = -4823.615  # was -4976.4165

$ws.Range("H32").Value = 490  # was 535
$ws.Range("I32").Value = 490  # was 535
$ws.Range("K32").Value = 490  # was 535
$ws.Range("M32").Value = -174  # was -219

$ws.Range("H34").Value = 4219.6787  # was 4402
$ws.Range("I34").Value = 4207.6  # was 4415.357
$ws.Range("J34").Value = 4233.615  # was 4386.4165
$ws.Range("K34").Value = 4207.6  # was 4415.357
$ws.Range("L34").Value = 4233.615  # was 4386.4165
$ws.Range("M34").Value = -4005.6  # was -4213.357
$ws.Range("N34").Value = -4637.615  # was -4790.4165

$ws.Range("H58").Value = 4278.522  # was 4429.591
$ws.Range("I58").Value = 3592.6667  # was 3705.5715
$ws.Range("J58").Value = 5564.5  # was 5696.625
$ws.Range("K58").Value = 3592.6667  # was 3705.5715
$ws.Range("L58").Value = 5564.5  # was 5696.625
$ws.Range("M58").Value = -3389.6667  # was -3502.5715
$ws.Range("N58").Value = -5970.5  # was -6102.625

$ws.Range("H88").Value = 22295  # was 23876.691
$ws.Range("J88").Value = 21986.076  # was 23673.834
$ws.Range("L88").Value = 21986.076  # was 23673.834
$ws.Range("N88").Value = -22798.076  # was -24485.834

$ws.Range("H91").Value = 22295  # was 23876.691
$ws.Range("J91").Value = 21986.076  # was 23673.834
$ws.Range("L91").Value = 21986.076  # was 23673.834
$ws.Range("N91").Value = -24794.076  # was -26481.834

$ws.Range("H136").Value = 4278.522  # was 4429.591
$ws.Range("I136").Value = 3592.6667  # was 3705.5715
$ws.Range("J136").Value = 5564.5  # was 5696.625
$ws.Range("K136").Value = 10778.0001  # was 11116.7145
$ws.Range("L136").Value = 16693.5  # was 17089.875
$ws.Range("M136").Value = -8228.000100000001  # was -8566.7145
$ws.Range("N136").Value = -21793.5  # was -22189.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 41178924  # was 39758964
$ws.Range("I4").Value = 41178924  # was 39758964
$ws.Range("K4").Value = 123536772  # was 119276892
$ws.Range("M4").Value = -123536660  # was -119276780

$ws.Range("H40").Value = 70.90000000000001  # was 63.2
$ws.Range("I40").Value = 29.6  # was 25.666666
$ws.Range("J40").Value = 112.2  # was 119.5
$ws.Range("K40").Value = 118.4  # was 102.666664
$ws.Range("L40").Value = 448.8  # was 478
$ws.Range("M40").Value = -49.40000000000001  # was -33.666664
$ws.Range("N40").Value = -586.8  # was -616

$ws.Range("H139").Value = 1862.7  # was 1958.6666
$ws.Range("I139").Value = 1461  # was 1538
$ws.Range("K139").Value = 4383  # was 4614
$ws.Range("M139").Value = 757  # was 526

$ws.Range("H140").Value = 1666.3334  # was 2000
$ws.Range("I140").Value = 1666.3334  # was 2000
$ws.Range("K140").Value = 4999.0002  # was 6000
$ws.Range("M140").Value = 180.9997999999996  # was -820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3563.6316  # was 3610.45
$ws.Range("I97").Value = 2146.182  # was 2100.818
$ws.Range("J97").Value = 5512.625  # was 5455.5557
$ws.Range("K97").Value = 2146.182  # was 2100.818
$ws.Range("L97").Value = 5512.625  # was 5455.5557
$ws.Range("M97").Value = -1650.182  # was -1604.818
$ws.Range("N97").Value = -6504.625  # was -6447.5557

$ws.Range("H98").Value = 33029.855  # was 21401.75
$ws.Range("J98").Value = 33029.855  # was 21401.75
$ws.Range("L98").Value = 33029.855  # was 21401.75
$ws.Range("N98").Value = -39019.855  # was -27391.75

$ws.Range("H126").Value = 4094.25  # was 3876.0527
$ws.Range("I126").Value = 3550.8  # was 3513.2727
$ws.Range("J126").Value = 5000  # was 4374.875
$ws.Range("K126").Value = 10652.4  # was 10539.8181
$ws.Range("L126").Value = 15000  # was 13124.625
$ws.Range("M126").Value = -8182.400000000001  # was -8069.8181
$ws.Range("N126").Value = -19940  # was -18064.625

$ws.Range("H132").Value = 2481.9546  # was 2439.8696
$ws.Range("I132").Value = 1966  # was 2091
$ws.Range("J132").Value = 3101.1  # was 2820.4546
$ws.Range("K132").Value = 5898  # was 6273
$ws.Range("L132").Value = 9303.299999999999  # was 8461.363799999999
$ws.Range("M132").Value = -3368  # was -3743
$ws.Range("N132").Value = -14363.3  # was -13521.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5722.25  # was 4208.778
$ws.Range("I7").Value = 5796.3335  # was 4675.8
$ws.Range("J7").Value = 5500  # was 3625
$ws.Range("K7").Value = 5796.3335  # was 4675.8
$ws.Range("L7").Value = 5500  # was 3625
$ws.Range("M7").Value = -5684.3335  # was -4563.8
$ws.Range("N7").Value = -5724  # was -3849

$ws.Range("H16").Value = 49998  # was 26249
$ws.Range("J16").Value = 0  # was 2500
$ws.Range("L16").Value = 0  # was 2500
$ws.Range("N16").ClearContents()  # was -2840

$ws.Range("H40").Value = 7829.9165  # was 7503.846
$ws.Range("I40").Value = 7905.364  # was 7959.091
$ws.Range("J40").Value = 7000  # was 5000
$ws.Range("K40").Value = 7905.364  # was 7959.091
$ws.Range("L40").Value = 7000  # was 5000
$ws.Range("M40").Value = -7769.364  # was -7823.091
$ws.Range("N40").Value = -7272  # was -5272

$ws.Range("H47").Value = 0  # was 20000
$ws.Range("I47").Value = 0  # was 20000
$ws.Range("K47").Value = 0  # was 20000
$ws.Range("M47").ClearContents()  # was -19510

$ws.Range("H52").Value = 0  # was 20000
$ws.Range("I52").Value = 0  # was 20000
$ws.Range("K52").Value = 0  # was 20000
$ws.Range("M52").ClearContents()  # was -19767

$ws.Range("H68").Value = 4633.3335  # was 3900
$ws.Range("I68").Value = 4450  # was 3350
$ws.Range("K68").Value = 4450  # was 3350
$ws.Range("M68").Value = -3701  # was -2601

$ws.Range("H71").Value = 4633.3335  # was 3900
$ws.Range("I71").Value = 4450  # was 3350
$ws.Range("K71").Value = 22250  # was 16750
$ws.Range("M71").Value = -18506  # was -13006

$ws.Range("H126").Value = 5722.25  # was 4208.778
$ws.Range("I126").Value = 5796.3335  # was 4675.8
$ws.Range("J126").Value = 5500  # was 3625
$ws.Range("K126").Value = 17389.0005  # was 14027.4
$ws.Range("L126").Value = 16500  # was 10875
$ws.Range("M126").Value = -14919.0005  # was -11557.4
$ws.Range("N126").Value = -21440  # was -15815

$ws.Range("H132").Value = 2950  # was 3166.3333
$ws.Range("I132").Value = 3000  # was 2249.5
$ws.Range("J132").Value = 2900  # was 5000
$ws.Range("K132").Value = 9000  # was 6748.5
$ws.Range("L132").Value = 8700  # was 15000
$ws.Range("M132").Value = -6470  # was -4218.5
$ws.Range("N132").Value = -13760  # was -20060

$ws.Range("H136").Value = 2714.9473  # was 2639.2
$ws.Range("J136").Value = 4997  # was 4047.75
$ws.Range("L136").Value = 14991  # was 12143.25
$ws.Range("N136").Value = -20091  # was -17243.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 41374.5  # was 41832.668
$ws.Range("J47").Value = 41374.5  # was 41832.668
$ws.Range("L47").Value = 41374.5  # was 41832.668
$ws.Range("N47").Value = -42518.5  # was -42976.668

$ws.Range("H132").Value = 3332.7222  # was 3112.3333
$ws.Range("I132").Value = 1765.1428  # was 1670.9
$ws.Range("J132").Value = 8819.25  # was 7917.1113
$ws.Range("K132").Value = 5295.428400000001  # was 5012.700000000001
$ws.Range("L132").Value = 26457.75  # was 23751.3339
$ws.Range("M132").Value = -2765.428400000001  # was -2482.700000000001
$ws.Range("N132").Value = -31517.75  # was -28811.3339

